$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: Model Cybertruck -> Model x, Year 0 -> 2019, Color Silver -> White
$ws.Range("C30").Value = "Model x"
$ws.Range("D30").Value = 2019
$ws.Range("G30").Value = "White"

# Row 31: Model Cybertruck -> Model x, Year 0 -> 2018, Color Silver -> Black
$ws.Range("C31").Value = "Model x"
$ws.Range("D31").Value = 2018
$ws.Range("G31").Value = "Black"

# Row 34: Engine Type Standard -> Electric
$ws.Range("E34").Value = "Electric"

# Row 35: Engine Type Mid Range -> Electric, Color Red -> Blue
$ws.Range("E35").Value = "Electric"
$ws.Range("G35").Value = "Blue"

# Row 36: Engine Type Performance -> Electric, Color Red -> White
$ws.Range("E36").Value = "Electric"
$ws.Range("G36").Value = "White"
# E36 reverts from the group's "last row" thick-bottom border style to the regular thin-bottom style
$ws.Range("E36").Borders.Item(9).Weight = 2

# Update the current selection/view to match the saved workbook's last state
$ws.Range("G33").Select()
